$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits
#    right under the H1 title.
# ------------------------------------------------------------------
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Meta description*") {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -gt 0) {
    $d.Paragraphs.Item($metaIndex).Range.Delete()
}

# ------------------------------------------------------------------
# 2. Locate the trailing "Create a feature image for Aztec Wilds..."
#    (italic) paragraph at the end of the document.
# ------------------------------------------------------------------
$imgPromptIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Create a feature image for Aztec Wilds*") {
        $imgPromptIndex = $i
        break
    }
}

if ($imgPromptIndex -gt 0) {
    $imgPromptPara = $d.Paragraphs.Item($imgPromptIndex)

    # 2a. Insert a brand-new bold paragraph ("Play Aztec Wilds for Free -
    #     Review of Unique Gameplay and High RTP") right before it, matching
    #     the existing "leading empty run + bold run" pattern used elsewhere
    #     in the document.
    $imgPromptPara.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs.Item($imgPromptIndex)
    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aztec Wilds for Free - Review of Unique Gameplay and High RTP</w:t></w:r></w:p>'
    $newPara.Range.InsertXML($newParaXml)

    # Index of the image-prompt paragraph shifted down by one.
    $imgPromptIndex = $imgPromptIndex + 1
}

# ------------------------------------------------------------------
# 3. Replace the text of the (still italic) image-prompt paragraph with
#    the meta-description copy, keeping its existing italic formatting.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a feature image for Aztec Wilds that features a happy Maya warrior wearing glasses in cartoon style. The warrior should be depicted in vibrant colors and holding a golden treasure chest to represent the hidden treasures that players are aiming to uncover in the game. Add elements of the Aztec civilization in the background, such as pyramids and towering stone sculptures, to create an immersive and authentic feel. Use bold fonts to highlight the game name and add a tagline, such as `"Discover the Hidden Treasures of the Aztecs!`" to entice players to explore the game's unique features.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Experience cascading symbols and three types of Wilds in Aztec Wilds. High RTP and immersive music. Play for free.",
    2)
